$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Re-theme the deck: swap the slide master's colour scheme from the
#    "Integral" (Red Violet) palette to the stock "Office" palette.
#    (ThemeColorScheme.Colors(i).RGB expects a Win32 COLORREF, i.e. the
#    bytes of the RRGGBB hex value reversed to BBGGRR.)
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink 954F72

# ------------------------------------------------------------------
# 2) Re-style the three tables (slides 14-16) from the bespoke
#    "no style / no grid" table style to the built-in
#    "Medium Style 2 - Accent 1" table style.
# ------------------------------------------------------------------
$newTableStyle = "{5C136576-62D1-4677-9D37-C8096D437757}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
